$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 100800
$ws.Range("E2").Value = 327272510
$ws.Range("C3").Value = 249273
$ws.Range("E3").Value = 1036172372
$ws.Range("C5").Value = 39477
$ws.Range("E5").Value = 361411392
$ws.Range("C46").Value = 10920
$ws.Range("E46").Value = 66193511
$ws.Range("C53").Value = 141655
$ws.Range("E53").Value = 589981563
$ws.Range("C55").Value = 23187
$ws.Range("E55").Value = 187773972
$ws.Range("C57").Value = 3705
$ws.Range("E57").Value = 137905354
$ws.Range("C63").Value = 14101
$ws.Range("E63").Value = 35459928
$ws.Range("C64").Value = 5055
$ws.Range("E64").Value = 19389984
$ws.Range("C79").Value = 116574
$ws.Range("E79").Value = 447276129
$ws.Range("C81").Value = 17423
$ws.Range("E81").Value = 133403196
$ws.Range("C91").Value = 150984
$ws.Range("E91").Value = 480931655
$ws.Range("C92").Value = 408705
$ws.Range("E92").Value = 1590402401
$ws.Range("C93").Value = 209256
$ws.Range("E93").Value = 1303951700
$ws.Range("C94").Value = 93999
$ws.Range("E94").Value = 911298136
$ws.Range("C96").Value = 17144
$ws.Range("E96").Value = 780902892
$ws.Range("C104").Value = 135152
$ws.Range("E104").Value = 271713428
$ws.Range("C106").Value = 18118
$ws.Range("E106").Value = 40712935
$ws.Range("C114").Value = 3714
$ws.Range("E114").Value = 8915019
$ws.Range("C115").Value = 11468
$ws.Range("E115").Value = 32283158
$ws.Range("C116").Value = 4424
$ws.Range("E116").Value = 19530976
$ws.Range("C118").Value = 906
$ws.Range("E118").Value = 10421491
$ws.Range("C122").Value = 8324
$ws.Range("E122").Value = 12555568
$ws.Range("C131").Value = 75573
$ws.Range("E131").Value = 307076240
$ws.Range("C166").Value = 35925
$ws.Range("E166").Value = 210542620
$ws.Range("C174").Value = 226046
$ws.Range("E174").Value = 900132284
$ws.Range("C175").Value = 80749
$ws.Range("E175").Value = 485283850
$ws.Range("C184").Value = 68726
$ws.Range("E184").Value = 134113858
